$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 35719510
$ws.Range("J70").Value = 50006216
$ws.Range("L70").Value = 150018648
$ws.Range("N70").Value = -150019188

$ws.Range("H73").Value = 35719510
$ws.Range("J73").Value = 50006216
$ws.Range("L73").Value = 150018648
$ws.Range("N73").Value = -150020520

$ws.Range("H86").Value = 100005224
$ws.Range("I86").Value = 5702.8
$ws.Range("K86").Value = 5702.8
$ws.Range("M86").Value = -4579.8

$ws.Range("H89").Value = 100005224
$ws.Range("I89").Value = 5702.8
$ws.Range("K89").Value = 28514
$ws.Range("M89").Value = -22898

$ws.Range("H100").Value = 3261.7693
$ws.Range("I100").Value = 1817.1666
$ws.Range("J100").Value = 4500
$ws.Range("K100").Value = 1817.1666
$ws.Range("L100").Value = 4500
$ws.Range("M100").Value = -1276.1666
$ws.Range("N100").Value = -5582

$ws.Range("H107").Value = 3602.6538
$ws.Range("I107").Value = 3107.2632
$ws.Range("J107").Value = 4947.2856
$ws.Range("K107").Value = 3107.2632
$ws.Range("L107").Value = 4947.2856
$ws.Range("M107").Value = -1187.2632
$ws.Range("N107").Value = -8787.285599999999

$ws.Range("H113").Value = 2351.4285
$ws.Range("I113").Value = 2584.3333
$ws.Range("J113").Value = 2176.75
$ws.Range("K113").Value = 2584.3333
$ws.Range("L113").Value = 2176.75
$ws.Range("M113").Value = 669.6667000000002
$ws.Range("N113").Value = -8684.75

$ws.Range("H117").Value = 110740.6
$ws.Range("J117").Value = 110740.6
$ws.Range("L117").Value = 110740.6
$ws.Range("N117").Value = -119918.6

$ws.Range("H125").Value = 1171
$ws.Range("I125").Value = 1333.6666
$ws.Range("J125").Value = 1049
$ws.Range("K125").Value = 12002.9994
$ws.Range("L125").Value = 9441
$ws.Range("M125").Value = -9542.999400000001
$ws.Range("N125").Value = -14361

$ws.Range("H137").Value = 4979.1763
$ws.Range("I137").Value = 3850.25
$ws.Range("K137").Value = 11550.75
$ws.Range("M137").Value = -9000.75

$ws.Range("H138").Value = 3195.8372
$ws.Range("I138").Value = 2513.1875
$ws.Range("K138").Value = 7539.5625
$ws.Range("M138").Value = -2399.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10968017
$ws.Range("I32").Value = 6174321
$ws.Range("K32").Value = 6174321
$ws.Range("M32").Value = -6174034

$ws.Range("H62").Value = 54997.5
$ws.Range("J62").Value = 54997.5
$ws.Range("L62").Value = 54997.5
$ws.Range("N62").Value = -56245.5

$ws.Range("H65").Value = 54997.5
$ws.Range("J65").Value = 54997.5
$ws.Range("L65").Value = 164992.5
$ws.Range("N65").Value = -171232.5

$ws.Range("H102").Value = 2766.5
$ws.Range("I102").Value = 2309.4285
$ws.Range("J102").Value = 3833
$ws.Range("K102").Value = 2309.4285
$ws.Range("L102").Value = 3833
$ws.Range("M102").Value = -687.4285
$ws.Range("N102").Value = -7077

$ws.Range("H110").Value = 2124.5
$ws.Range("I110").Value = 1311.75
$ws.Range("J110").Value = 3750
$ws.Range("K110").Value = 1311.75
$ws.Range("L110").Value = 3750
$ws.Range("M110").Value = 733.25
$ws.Range("N110").Value = -7840

$ws.Range("H132").Value = 2665.4792
$ws.Range("I132").Value = 2307.838
$ws.Range("K132").Value = 6923.514000000001
$ws.Range("M132").Value = -4393.514000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2986.8572
$ws.Range("I99").Value = 2901.3333
$ws.Range("K99").Value = 2901.3333
$ws.Range("M99").Value = -1403.3333

$ws.Range("H105").Value = 3461.2104
$ws.Range("I105").Value = 2779.889
$ws.Range("K105").Value = 2779.889
$ws.Range("M105").Value = -1032.889

$ws.Range("H132").Value = 117680.75
$ws.Range("J132").Value = 117680.75
$ws.Range("L132").Value = 117680.75
$ws.Range("N132").Value = -127800.75

$ws.Range("H134").Value = 11566819
$ws.Range("I134").Value = 2977775.5
$ws.Range("J134").Value = 30306552
$ws.Range("K134").Value = 8933326.5
$ws.Range("L134").Value = 90919656
$ws.Range("M134").Value = -8930791.5
$ws.Range("N134").Value = -90924726

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5144.4736
$ws.Range("I31").Value = 3146.3333
$ws.Range("K31").Value = 3146.3333
$ws.Range("M31").Value = -2851.3333

$ws.Range("H34").Value = 5144.4736
$ws.Range("I34").Value = 3146.3333
$ws.Range("K34").Value = 3146.3333
$ws.Range("M34").Value = -2944.3333

$ws.Range("H132").Value = 2249.1667
$ws.Range("I132").Value = 873.75
$ws.Range("K132").Value = 2621.25
$ws.Range("M132").Value = -91.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 26.181818
$ws.Range("I2").Value = 8.166667
$ws.Range("K2").Value = 49.000002
$ws.Range("M2").Value = 63.999998

$ws.Range("H68").Value = 1595.2727
$ws.Range("I68").Value = 1582.3334
$ws.Range("J68").Value = 1610.8
$ws.Range("K68").Value = 4747.0002
$ws.Range("L68").Value = 4832.4
$ws.Range("M68").Value = -3936.0002
$ws.Range("N68").Value = -6454.4

$ws.Range("H71").Value = 1595.2727
$ws.Range("I71").Value = 1582.3334
$ws.Range("J71").Value = 1610.8
$ws.Range("K71").Value = 14241.0006
$ws.Range("L71").Value = 14497.2
$ws.Range("M71").Value = -10185.0006
$ws.Range("N71").Value = -22609.2

$ws.Range("H131").Value = 1686.2363
$ws.Range("J131").Value = 1855.2927
$ws.Range("L131").Value = 5565.8781
$ws.Range("N131").Value = -15645.8781

$ws.Range("H132").Value = 1698.4375
$ws.Range("J132").Value = 1924.3334
$ws.Range("L132").Value = 17319.0006
$ws.Range("N132").Value = -22379.0006

$ws.Range("H140").Value = 18184400
$ws.Range("I140").Value = 18184400
$ws.Range("K140").Value = 54553200
$ws.Range("M140").Value = -54548020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2421.5
$ws.Range("I31").Value = 2421.5
$ws.Range("K31").Value = 2421.5
$ws.Range("M31").Value = -2129.5

$ws.Range("H37").Value = 2421.5
$ws.Range("I37").Value = 2421.5
$ws.Range("K37").Value = 2421.5
$ws.Range("M37").Value = -2144.5

$ws.Range("H107").Value = 932.9231
$ws.Range("I107").Value = 963
$ws.Range("J107").Value = 897.8333
$ws.Range("K107").Value = 963
$ws.Range("L107").Value = 897.8333
$ws.Range("M107").Value = 957
$ws.Range("N107").Value = -4737.8333

$ws.Range("H132").Value = 1707.4117
$ws.Range("I132").Value = 1334.1333
$ws.Range("J132").Value = 4507
$ws.Range("K132").Value = 4002.3999
$ws.Range("L132").Value = 13521
$ws.Range("M132").Value = -1472.3999
$ws.Range("N132").Value = -18581

$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 269.12903
$ws.Range("I55").Value = 206.75
$ws.Range("J55").Value = 335.66666
$ws.Range("K55").Value = 206.75
$ws.Range("L55").Value = 335.66666
$ws.Range("M55").Value = -33.75
$ws.Range("N55").Value = -681.66666

$ws.Range("H68").Value = 3291.4285
$ws.Range("I68").Value = 2020.5
$ws.Range("K68").Value = 2020.5
$ws.Range("M68").Value = -1271.5

$ws.Range("H71").Value = 3291.4285
$ws.Range("I71").Value = 2020.5
$ws.Range("K71").Value = 10102.5
$ws.Range("M71").Value = -6358.5

$ws.Range("H82").Value = 2318.9614
$ws.Range("I82").Value = 1844.7693
$ws.Range("K82").Value = 1844.7693
$ws.Range("M82").Value = -1483.7693

$ws.Range("H85").Value = 2318.9614
$ws.Range("I85").Value = 1844.7693
$ws.Range("K85").Value = 1844.7693
$ws.Range("M85").Value = -596.7692999999999

$ws.Range("H93").Value = 1559.25
$ws.Range("I93").Value = 1246.5
$ws.Range("J93").Value = 1872
$ws.Range("K93").Value = 1246.5
$ws.Range("L93").Value = 1872
$ws.Range("M93").Value = 1.5
$ws.Range("N93").Value = -4368

$ws.Range("H99").Value = 74160.664
$ws.Range("J99").Value = 103323
$ws.Range("L99").Value = 103323
$ws.Range("N99").Value = -109313

$ws.Range("H100").Value = 2728.2856
$ws.Range("I100").Value = 2016.3334
$ws.Range("J100").Value = 7000
$ws.Range("K100").Value = 2016.3334
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -1475.3334
$ws.Range("N100").Value = -8082

$ws.Range("H132").Value = 4528.0713
$ws.Range("I132").Value = 4126.727
$ws.Range("K132").Value = 12380.181
$ws.Range("M132").Value = -9850.181

$ws.Range("H136").Value = 4629.4814
$ws.Range("I136").Value = 4061.4375
$ws.Range("K136").Value = 12184.3125
$ws.Range("M136").Value = -9634.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 41583.332
$ws.Range("J86").Value = 41583.332
$ws.Range("L86").Value = 41583.332
$ws.Range("N86").Value = -43829.332

$ws.Range("H89").Value = 41583.332
$ws.Range("J89").Value = 41583.332
$ws.Range("L89").Value = 207916.66
$ws.Range("N89").Value = -219148.66

$ws.Range("H92").Value = 44999.5
$ws.Range("I92").Value = 49999
$ws.Range("J92").Value = 40000
$ws.Range("K92").Value = 49999
$ws.Range("L92").Value = 40000
$ws.Range("M92").Value = -47503
$ws.Range("N92").Value = -44992

$ws.Range("H96").Value = 35677.855
$ws.Range("J96").Value = 45187.5
$ws.Range("L96").Value = 45187.5
$ws.Range("N96").Value = -47933.5
